$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Notes" column header (new shared string, dimension grows to C1)
$ws.Range("C1").Value = "Notes"

# Left-align the header row (new cellXfs entry: fontId 0 + alignment horizontal=left)
$ws.Range("A1:C1").HorizontalAlignment = -4131

# Resize the columns (close to the widths Excel recalculates for the new layout,
# column C is sized to fit the new "Notes" header)
$ws.Columns.Item(1).ColumnWidth = 14.571428571428571
$ws.Columns.Item(2).ColumnWidth = 17.571428571428571
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(3).ColumnWidth = 9.857142857142858

# Extend / move the live selection used in the saved view
$ws.Range("A2:XFD10").Select()
